$wb = $excel.ActiveWorkbook

# ALC!row19
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1206.9
$ws.Range("J19").Value = 978.4286
$ws.Range("L19").Value = 978.4286
$ws.Range("N19").Value = -1328.4286

# ALC!row40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1432.2222
$ws.Range("I40").Value = 1540
$ws.Range("J40").Value = 1297.5
$ws.Range("K40").Value = 1540
$ws.Range("L40").Value = 1297.5
$ws.Range("M40").Value = -1365
$ws.Range("N40").Value = -1647.5

# ALC!row64
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 41309.69
$ws.Range("I64").Value = 74017.64
$ws.Range("J64").Value = 3150.4167
$ws.Range("K64").Value = 74017.64
$ws.Range("L64").Value = 3150.4167
$ws.Range("M64").Value = -73769.64
$ws.Range("N64").Value = -3646.4167

# ALC!row67
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 41309.69
$ws.Range("I67").Value = 74017.64
$ws.Range("J67").Value = 3150.4167
$ws.Range("K67").Value = 74017.64
$ws.Range("L67").Value = 3150.4167
$ws.Range("M67").Value = -73159.64
$ws.Range("N67").Value = -4866.4167

# ALC!row121
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 1448
$ws.Range("J121").Value = 1562.6666
$ws.Range("L121").Value = 4687.9998
$ws.Range("N121").Value = -8181.9998

# ARM!row32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 29850.223
$ws.Range("I32").Value = 6673.4917
$ws.Range("J32").Value = 371707
$ws.Range("K32").Value = 6673.4917
$ws.Range("L32").Value = 371707
$ws.Range("M32").Value = -6386.4917
$ws.Range("N32").Value = -372281

# ARM!row45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 72707
$ws.Range("I45").Value = 92226.09
$ws.Range("J45").Value = 1137
$ws.Range("K45").Value = 92226.09
$ws.Range("L45").Value = 1137
$ws.Range("M45").Value = -91849.09
$ws.Range("N45").Value = -1891

# ARM!row74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1054.1538
$ws.Range("I74").Value = 1054.1538
$ws.Range("K74").Value = 1054.1538
$ws.Range("M74").Value = -180.1538

# ARM!row77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1054.1538
$ws.Range("I77").Value = 1054.1538
$ws.Range("K77").Value = 5270.769
$ws.Range("M77").Value = -902.7690000000002

# ARM!row110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 40044796
$ws.Range("I110").Value = 47672204
$ws.Range("K110").Value = 47672204
$ws.Range("M110").Value = -47670159

# ARM!row122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1670.9333
$ws.Range("I122").Value = 1561.4286
$ws.Range("J122").Value = 1926.4445
$ws.Range("K122").Value = 4684.2858
$ws.Range("L122").Value = 5779.333500000001
$ws.Range("M122").Value = -2234.2858
$ws.Range("N122").Value = -10679.3335

# ARM!row132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4626.174
$ws.Range("I132").Value = 4959.1177
$ws.Range("J132").Value = 3682.8333
$ws.Range("K132").Value = 14877.3531
$ws.Range("L132").Value = 11048.4999
$ws.Range("M132").Value = -12347.3531
$ws.Range("N132").Value = -16108.4999

# CRP!row58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2105.7058
$ws.Range("I58").Value = 1600.8572
$ws.Range("J58").Value = 2459.1
$ws.Range("K58").Value = 1600.8572
$ws.Range("L58").Value = 2459.1
$ws.Range("M58").Value = -1397.8572
$ws.Range("N58").Value = -2865.1

# CRP!row136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2105.7058
$ws.Range("I136").Value = 1600.8572
$ws.Range("J136").Value = 2459.1
$ws.Range("K136").Value = 4802.571599999999
$ws.Range("L136").Value = 7377.299999999999
$ws.Range("M136").Value = -2252.571599999999
$ws.Range("N136").Value = -12477.3

# CUL!row68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1400.6666
$ws.Range("I68").Value = 1100
$ws.Range("J68").Value = 1486.5714
$ws.Range("K68").Value = 3300
$ws.Range("L68").Value = 4459.7142
$ws.Range("M68").Value = -2489
$ws.Range("N68").Value = -6081.7142

# CUL!row71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 1400.6666
$ws.Range("I71").Value = 1100
$ws.Range("J71").Value = 1486.5714
$ws.Range("K71").Value = 9900
$ws.Range("L71").Value = 13379.1426
$ws.Range("M71").Value = -5844
$ws.Range("N71").Value = -21491.1426

# CUL!row98
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 111458.555
$ws.Range("J98").Value = 166907.83
$ws.Range("L98").Value = 500723.49
$ws.Range("N98").Value = -503719.49

# CUL!row120
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 338343.34
$ws.Range("I120").Value = 338343.34
$ws.Range("K120").Value = 1015030.02
$ws.Range("M120").Value = -1010192.02

# CUL!row131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 794.19385
$ws.Range("I131").Value = 224
$ws.Range("J131").Value = 824.8495
$ws.Range("K131").Value = 672
$ws.Range("L131").Value = 2474.5485
$ws.Range("M131").Value = 4368
$ws.Range("N131").Value = -12554.5485

# CUL!row132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1601.5483
$ws.Range("I132").Value = 794.8333
$ws.Range("J132").Value = 2111.0527
$ws.Range("K132").Value = 7153.4997
$ws.Range("L132").Value = 18999.4743
$ws.Range("M132").Value = -4623.4997
$ws.Range("N132").Value = -24059.4743

# CUL!row139
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 1731.7407
$ws.Range("I139").Value = 767.58826
$ws.Range("J139").Value = 3370.8
$ws.Range("K139").Value = 2302.76478
$ws.Range("L139").Value = 10112.4
$ws.Range("M139").Value = 2837.23522
$ws.Range("N139").Value = -20392.4

# CUL!row140
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 5572
$ws.Range("I140").Value = 6524.6
$ws.Range("J140").Value = 2396.6667
$ws.Range("K140").Value = 19573.8
$ws.Range("L140").Value = 7190.000100000001
$ws.Range("M140").Value = -14393.8
$ws.Range("N140").Value = -17550.0001

# GSM!row97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 111113624
$ws.Range("I97").Value = 111113624
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 111113624
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -111113128
$ws.Range("N97").ClearContents()

# GSM!row122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 782.8889
$ws.Range("I122").Value = 488.2353
$ws.Range("J122").Value = 1283.8
$ws.Range("K122").Value = 1464.7059
$ws.Range("L122").Value = 3851.4
$ws.Range("M122").Value = 985.2941000000001
$ws.Range("N122").Value = -8751.4

# LTW!row7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3408.6667
$ws.Range("I7").Value = 2743.4285
$ws.Range("J7").Value = 4340
$ws.Range("K7").Value = 2743.4285
$ws.Range("L7").Value = 4340
$ws.Range("M7").Value = -2631.4285
$ws.Range("N7").Value = -4564

# LTW!row126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 3408.6667
$ws.Range("I126").Value = 2743.4285
$ws.Range("J126").Value = 4340
$ws.Range("K126").Value = 8230.2855
$ws.Range("L126").Value = 13020
$ws.Range("M126").Value = -5760.2855
$ws.Range("N126").Value = -17960

# LTW!row136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2011.6538
$ws.Range("I136").Value = 1834.9131
$ws.Range("J136").Value = 3366.6667
$ws.Range("K136").Value = 5504.7393
$ws.Range("L136").Value = 10100.0001
$ws.Range("M136").Value = -2954.7393
$ws.Range("N136").Value = -15200.0001

# WVR!row109
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

# WVR!row126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1776.75
$ws.Range("I126").Value = 1459.1428
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 4377.428400000001
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -1907.428400000001
$ws.Range("N126").Value = -16940

# WVR!row132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1863.4259
$ws.Range("I132").Value = 1824.3572
$ws.Range("J132").Value = 2000.1666
$ws.Range("K132").Value = 5473.071599999999
$ws.Range("L132").Value = 6000.4998
$ws.Range("M132").Value = -2943.071599999999
$ws.Range("N132").Value = -11060.4998

# WVR!row136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1556.9166
$ws.Range("I136").Value = 668.9583
$ws.Range("J136").Value = 3332.8333
$ws.Range("K136").Value = 2006.8749
$ws.Range("L136").Value = 9998.499899999999
$ws.Range("M136").Value = 543.1251
$ws.Range("N136").Value = -15098.4999

Write-Host "Update complete"